$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.032.67'
$ws.Range('D2').Style = $ws.Range('B2').Style
$ws.Range('E2').Value = '  +7.57%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.875.30'
$ws.Range('D3').Style = $ws.Range('B3').Style
$ws.Range('E3').Value = '  +5.46%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = $ws.Range('B4').Style

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.90'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  +2.30%  '

$ws.Range('E6').Value = '  +0.15%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4974'
$ws.Range('D7').Style = $ws.Range('B7').Style
$ws.Range('E7').Value = '  +1.68%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.70'
$ws.Range('D8').Style = $ws.Range('B8').Style
$ws.Range('E8').Value = '  +9.16%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2843'
$ws.Range('D9').Style = $ws.Range('B9').Style
$ws.Range('E9').Value = '  +7.36%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06554'
$ws.Range('D10').Style = $ws.Range('B10').Style
$ws.Range('E10').Value = '  +5.12%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.872.53'
$ws.Range('D11').Style = $ws.Range('B11').Style
$ws.Range('E11').Value = '  +5.27%  '

$ws.Range('E12').Value = '  +4.93%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07187'
$ws.Range('D13').Style = $ws.Range('B13').Style
$ws.Range('E13').Value = '  +2.63%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6606'
$ws.Range('D14').Style = $ws.Range('B14').Style
$ws.Range('E14').Value = '  +7.37%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '85.09'
$ws.Range('D15').Style = $ws.Range('B15').Style
$ws.Range('E15').Value = '  +7.28%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.799'
$ws.Range('D16').Style = $ws.Range('B16').Style
$ws.Range('E16').Value = '  +4.29%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.022.44'
$ws.Range('D17').Style = $ws.Range('B17').Style
$ws.Range('E17').Value = '  +7.64%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9987'
$ws.Range('D18').Style = $ws.Range('B18').Style
$ws.Range('E18').Value = '  -0.19%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.87'
$ws.Range('D19').Style = $ws.Range('B19').Style
$ws.Range('E19').Value = '  +9.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007509'
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  +4.35%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  +0.08%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.112.13'
$ws.Range('D22').Style = $ws.Range('B22').Style
$ws.Range('E22').Value = '  +5.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.745'
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  +4.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.021'
$ws.Range('D24').Style = $ws.Range('B24').Style
$ws.Range('E24').Value = '  +4.61%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.502'
$ws.Range('D25').Style = $ws.Range('B25').Style
$ws.Range('E25').Value = '  +5.92%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '144.41'
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  +1.86%  '

$ws.Range('E27').Value = '  +23.53%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.72'
$ws.Range('D28').Style = $ws.Range('B28').Style
$ws.Range('E28').Value = '  +7.39%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.955'
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  +5.32%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.386'
$ws.Range('D30').Style = $ws.Range('B30').Style
$ws.Range('E30').Value = '  -0.54%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.207'
$ws.Range('D31').Style = $ws.Range('B31').Style
$ws.Range('E31').Value = '  +2.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08595'
$ws.Range('D32').Style = $ws.Range('B32').Style
$ws.Range('E32').Value = '  +4.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.882'
$ws.Range('D33').Style = $ws.Range('B33').Style
$ws.Range('E33').Value = '  +2.79%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05068'
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  +6.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.138'
$ws.Range('D35').Style = $ws.Range('B35').Style
$ws.Range('E35').Value = '  +7.78%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').Style = $ws.Range('B36').Style
$ws.Range('E36').Value = '  +0.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6843'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  +6.75%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.704'
$ws.Range('D38').Style = $ws.Range('B38').Style
$ws.Range('E38').Value = '  +4.21%  '

$ws.Range('E39').Value = '  +13.67%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.741'
$ws.Range('D40').Style = $ws.Range('B40').Style
$ws.Range('E40').Value = '  +5.97%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9620'
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  +2.18%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01626'
$ws.Range('D42').Style = $ws.Range('B42').Style
$ws.Range('E42').Value = '  +6.01%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.073'
$ws.Range('D43').Style = $ws.Range('B43').Style
$ws.Range('E43').Value = '  +3.10%  '

$ws.Range('E44').Value = '  +0.00%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.09'
$ws.Range('D45').Style = $ws.Range('B45').Style
$ws.Range('E45').Value = '  +2.93%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4187'
$ws.Range('D46').Style = $ws.Range('B46').Style
$ws.Range('E46').Value = '  +6.44%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.498'
$ws.Range('D47').Style = $ws.Range('B47').Style
$ws.Range('E47').Value = '  +4.66%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1253'
$ws.Range('D48').Style = $ws.Range('B48').Style
$ws.Range('E48').Value = '  +5.11%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05623'
$ws.Range('D49').Style = $ws.Range('B49').Style
$ws.Range('E49').Value = '  +3.95%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '32.44'
$ws.Range('D50').Style = $ws.Range('B50').Style
$ws.Range('E50').Value = '  +6.87%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.205'
$ws.Range('D51').Style = $ws.Range('B51').Style
$ws.Range('E51').Value = '  +2.84%  '

